$d = $word.ActiveDocument

# 1. Merge "Android specificati" + "o" + "n link" into one run of text
$d.Content.Find.Execute("Android specification link", $true, $false, $false, $false, $false, $true, 1, $false, "Android specification link", 2) | Out-Null

# 3. Replace hyperlink display text (URL) with "PlayMarket: Netherlands Town Hall"
$d.Content.Find.Execute("https://play.google.com/store/apps/details?id=com.blackstone.ratusha", $true, $false, $false, $false, $false, $true, 1, $false, "PlayMarket: Netherlands Town Hall", 2) | Out-Null
